$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (76 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2303.7827
$ws.Range("I28").Value = 1923.8823
$ws.Range("J28").Value = 3380.1667
$ws.Range("K28").Value = 1923.8823
$ws.Range("L28").Value = 3380.1667
$ws.Range("M28").Value = -1438.8823
$ws.Range("N28").Value = -4350.1667
$ws.Range("H62").Value = 4842.4287
$ws.Range("I62").Value = 3699.5
$ws.Range("J62").Value = 5299.6
$ws.Range("K62").Value = 3699.5
$ws.Range("L62").Value = 5299.6
$ws.Range("M62").Value = -3075.5
$ws.Range("N62").Value = -6547.6
$ws.Range("H65").Value = 4842.4287
$ws.Range("I65").Value = 3699.5
$ws.Range("J65").Value = 5299.6
$ws.Range("K65").Value = 18497.5
$ws.Range("L65").Value = 26498
$ws.Range("M65").Value = -15377.5
$ws.Range("N65").Value = -32738
$ws.Range("H70").Value = 50006800
$ws.Range("I70").Value = 7900
$ws.Range("J70").Value = 55562230
$ws.Range("K70").Value = 23700
$ws.Range("L70").Value = 166686690
$ws.Range("M70").Value = -23430
$ws.Range("N70").Value = -166687230
$ws.Range("H73").Value = 50006800
$ws.Range("I73").Value = 7900
$ws.Range("J73").Value = 55562230
$ws.Range("K73").Value = 23700
$ws.Range("L73").Value = 166686690
$ws.Range("M73").Value = -22764
$ws.Range("N73").Value = -166688562
$ws.Range("H107").Value = 555.25
$ws.Range("I107").Value = 555.25
$ws.Range("K107").Value = 555.25
$ws.Range("M107").Value = 1364.75
$ws.Range("H112").Value = 2813.3281
$ws.Range("J112").Value = 2825.459
$ws.Range("L112").Value = 8476.377
$ws.Range("N112").Value = -10692.377
$ws.Range("H121").Value = 3654.9614
$ws.Range("J121").Value = 3654.9614
$ws.Range("L121").Value = 10964.8842
$ws.Range("N121").Value = -14458.8842
$ws.Range("H129").Value = 1450.5454
$ws.Range("I129").Value = 1064.3334
$ws.Range("K129").Value = 3193.0002
$ws.Range("M129").Value = 1806.9998
$ws.Range("H135").Value = 3288.366
$ws.Range("I135").Value = 797.129
$ws.Range("J135").Value = 11011.2
$ws.Range("K135").Value = 7174.161
$ws.Range("L135").Value = 99100.8
$ws.Range("M135").Value = -4639.161
$ws.Range("N135").Value = -104170.8
$ws.Range("H137").Value = 14495232
$ws.Range("I137").Value = 2620.4167
$ws.Range("J137").Value = 30305354
$ws.Range("K137").Value = 7861.250100000001
$ws.Range("L137").Value = 90916062
$ws.Range("M137").Value = -5311.250100000001
$ws.Range("N137").Value = -90921162
$ws.Range("H138").Value = 4212.3877
$ws.Range("I138").Value = 1567.4762
$ws.Range("J138").Value = 6196.0713
$ws.Range("K138").Value = 4702.4286
$ws.Range("L138").Value = 18588.2139
$ws.Range("M138").Value = 437.5713999999998
$ws.Range("N138").Value = -28868.2139
$ws.Range("H141").Value = 3485.1177
$ws.Range("I141").Value = 3348.303
$ws.Range("K141").Value = 10044.909
$ws.Range("M141").Value = -4864.909

# ---- Sheet: ARM (60 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9800.878000000001
$ws.Range("I32").Value = 9379.402
$ws.Range("J32").Value = 13134.363
$ws.Range("K32").Value = 9379.402
$ws.Range("L32").Value = 13134.363
$ws.Range("M32").Value = -9092.402
$ws.Range("N32").Value = -13708.363
$ws.Range("H43").Value = 22320
$ws.Range("J43").Value = 22320
$ws.Range("L43").Value = 22320
$ws.Range("N43").Value = -22946
$ws.Range("H61").Value = 5379.8423
$ws.Range("I61").Value = 5132.077
$ws.Range("K61").Value = 5132.077
$ws.Range("M61").Value = -4920.077
$ws.Range("H74").Value = 1280.8286
$ws.Range("I74").Value = 896.7917
$ws.Range("J74").Value = 2118.7273
$ws.Range("K74").Value = 896.7917
$ws.Range("L74").Value = 2118.7273
$ws.Range("M74").Value = -22.79169999999999
$ws.Range("N74").Value = -3866.7273
$ws.Range("H77").Value = 1280.8286
$ws.Range("I77").Value = 896.7917
$ws.Range("J77").Value = 2118.7273
$ws.Range("K77").Value = 4483.9585
$ws.Range("L77").Value = 10593.6365
$ws.Range("M77").Value = -115.9584999999997
$ws.Range("N77").Value = -19329.6365
$ws.Range("H80").Value = 84999
$ws.Range("J80").Value = 84999
$ws.Range("L80").Value = 84999
$ws.Range("N80").Value = -86995
$ws.Range("H83").Value = 84999
$ws.Range("J83").Value = 84999
$ws.Range("L83").Value = 254997
$ws.Range("N83").Value = -264981
$ws.Range("H132").Value = 27136.889
$ws.Range("I132").Value = 42120.145
$ws.Range("J132").Value = 11001.077
$ws.Range("K132").Value = 126360.435
$ws.Range("L132").Value = 33003.231
$ws.Range("M132").Value = -123830.435
$ws.Range("N132").Value = -38063.231
$ws.Range("H136").Value = 5379.8423
$ws.Range("I136").Value = 5132.077
$ws.Range("K136").Value = 15396.231
$ws.Range("M136").Value = -12846.231
$ws.Range("H138").Value = 79371.60000000001
$ws.Range("J138").Value = 79371.60000000001
$ws.Range("L138").Value = 79371.60000000001
$ws.Range("N138").Value = -89651.60000000001
$ws.Range("H139").Value = 78715
$ws.Range("J139").Value = 78715
$ws.Range("L139").Value = 78715
$ws.Range("N139").Value = -88995
$ws.Range("H140").Value = 97196.28999999999
$ws.Range("J140").Value = 100895.664
$ws.Range("L140").Value = 100895.664
$ws.Range("N140").Value = -111255.664

# ---- Sheet: BSM (49 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 9999
$ws.Range("J19").Value = 9999
$ws.Range("L19").Value = 9999
$ws.Range("N19").Value = -10345
$ws.Range("H20").Value = 7764.8945
$ws.Range("J20").Value = 6100.8887
$ws.Range("L20").Value = 6100.8887
$ws.Range("N20").Value = -6594.8887
$ws.Range("H80").Value = 682.7857
$ws.Range("J80").Value = 777.1
$ws.Range("L80").Value = 777.1
$ws.Range("N80").Value = -2773.1
$ws.Range("H83").Value = 682.7857
$ws.Range("J83").Value = 777.1
$ws.Range("L83").Value = 3885.5
$ws.Range("N83").Value = -13869.5
$ws.Range("H88").Value = 64984
$ws.Range("J88").Value = 64984
$ws.Range("L88").Value = 64984
$ws.Range("N88").Value = -65796
$ws.Range("H91").Value = 64984
$ws.Range("J91").Value = 64984
$ws.Range("L91").Value = 64984
$ws.Range("N91").Value = -67792
$ws.Range("H94").Value = 442829.38
$ws.Range("I94").Value = 548590.9399999999
$ws.Range("J94").Value = 2156.3333
$ws.Range("K94").Value = 548590.9399999999
$ws.Range("L94").Value = 2156.3333
$ws.Range("M94").Value = -548139.9399999999
$ws.Range("N94").Value = -3058.3333
$ws.Range("H99").Value = 53729740
$ws.Range("I99").Value = 68057010
$ws.Range("K99").Value = 68057010
$ws.Range("M99").Value = -68055512
$ws.Range("H105").Value = 42859416
$ws.Range("I105").Value = 45456710
$ws.Range("J105").Value = 3999.5
$ws.Range("K105").Value = 45456710
$ws.Range("L105").Value = 3999.5
$ws.Range("M105").Value = -45454963
$ws.Range("N105").Value = -7493.5
$ws.Range("H134").Value = 2066.7817
$ws.Range("I134").Value = 2095.2058
$ws.Range("J134").Value = 2020.762
$ws.Range("K134").Value = 6285.617400000001
$ws.Range("L134").Value = 6062.286
$ws.Range("M134").Value = -3750.617400000001
$ws.Range("N134").Value = -11132.286

# ---- Sheet: CRP (63 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20410786
$ws.Range("I31").Value = 26318002
$ws.Range("J31").Value = 4038.3635
$ws.Range("K31").Value = 26318002
$ws.Range("L31").Value = 4038.3635
$ws.Range("M31").Value = -26317707
$ws.Range("N31").Value = -4628.363499999999
$ws.Range("H34").Value = 20410786
$ws.Range("I34").Value = 26318002
$ws.Range("J34").Value = 4038.3635
$ws.Range("K34").Value = 26318002
$ws.Range("L34").Value = 4038.3635
$ws.Range("M34").Value = -26317800
$ws.Range("N34").Value = -4442.363499999999
$ws.Range("H74").Value = 79992.5
$ws.Range("J74").Value = 79992.5
$ws.Range("L74").Value = 79992.5
$ws.Range("N74").Value = -81740.5
$ws.Range("H77").Value = 79992.5
$ws.Range("J77").Value = 79992.5
$ws.Range("L77").Value = 239977.5
$ws.Range("N77").Value = -248713.5
$ws.Range("H86").Value = 12310.333
$ws.Range("I86").Value = 12608.667
$ws.Range("K86").Value = 12608.667
$ws.Range("M86").Value = -11485.667
$ws.Range("H88").Value = 52025.5
$ws.Range("J88").Value = 52025.5
$ws.Range("L88").Value = 52025.5
$ws.Range("N88").Value = -52837.5
$ws.Range("H89").Value = 12310.333
$ws.Range("I89").Value = 12608.667
$ws.Range("K89").Value = 63043.335
$ws.Range("M89").Value = -57427.335
$ws.Range("H91").Value = 52025.5
$ws.Range("J91").Value = 52025.5
$ws.Range("L91").Value = 52025.5
$ws.Range("N91").Value = -54833.5
$ws.Range("H124").Value = 48581.5
$ws.Range("J124").Value = 48581.5
$ws.Range("L124").Value = 48581.5
$ws.Range("N124").Value = -53491.5
$ws.Range("H129").Value = 52155.2
$ws.Range("I129").Value = 39999
$ws.Range("J129").Value = 100780
$ws.Range("K129").Value = 39999
$ws.Range("L129").Value = 100780
$ws.Range("M129").Value = -34999
$ws.Range("N129").Value = -110780
$ws.Range("H132").Value = 43026148
$ws.Range("I132").Value = 63504090
$ws.Range("J132").Value = 22475.3
$ws.Range("K132").Value = 190512270
$ws.Range("L132").Value = 67425.89999999999
$ws.Range("M132").Value = -190509740
$ws.Range("N132").Value = -72485.89999999999
$ws.Range("H134").Value = 1826.8889
$ws.Range("I134").Value = 1839.4615
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 5518.3845
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -2983.3845
$ws.Range("N134").Value = -9570

# ---- Sheet: CUL (229 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2429688.8
$ws.Range("I4").Value = 514090.03
$ws.Range("J4").Value = 11688416
$ws.Range("K4").Value = 1542270.09
$ws.Range("L4").Value = 35065248
$ws.Range("M4").Value = -1542158.09
$ws.Range("N4").Value = -35065472
$ws.Range("H11").Value = 63138676
$ws.Range("I11").Value = 100312.75
$ws.Range("J11").Value = 164000050
$ws.Range("K11").Value = 300938.25
$ws.Range("L11").Value = 492000150
$ws.Range("M11").Value = -300798.25
$ws.Range("N11").Value = -492000430
$ws.Range("H63").Value = 7678.4287
$ws.Range("J63").Value = 9749.799999999999
$ws.Range("L63").Value = 29249.4
$ws.Range("N63").Value = -30747.4
$ws.Range("H64").Value = 125012660
$ws.Range("I64").Value = 333334340
$ws.Range("J64").Value = 19653
$ws.Range("K64").Value = 1000003020
$ws.Range("L64").Value = 58959
$ws.Range("M64").Value = -1000002750
$ws.Range("N64").Value = -59499
$ws.Range("H66").Value = 7678.4287
$ws.Range("J66").Value = 9749.799999999999
$ws.Range("L66").Value = 87748.2
$ws.Range("N66").Value = -95236.2
$ws.Range("H67").Value = 125012660
$ws.Range("I67").Value = 333334340
$ws.Range("J67").Value = 19653
$ws.Range("K67").Value = 1000003020
$ws.Range("L67").Value = 58959
$ws.Range("M67").Value = -1000002084
$ws.Range("N67").Value = -60831
$ws.Range("H68").Value = 1663.6666
$ws.Range("J68").Value = 1997.25
$ws.Range("L68").Value = 5991.75
$ws.Range("N68").Value = -7613.75
$ws.Range("H69").Value = 3783.182
$ws.Range("I69").Value = 3721.25
$ws.Range("J69").Value = 3796.9443
$ws.Range("K69").Value = 11163.75
$ws.Range("L69").Value = 11390.8329
$ws.Range("M69").Value = -10352.75
$ws.Range("N69").Value = -13012.8329
$ws.Range("H70").Value = 166677330
$ws.Range("J70").Value = 166677330
$ws.Range("L70").Value = 500031990
$ws.Range("N70").Value = -500032620
$ws.Range("H71").Value = 1663.6666
$ws.Range("J71").Value = 1997.25
$ws.Range("L71").Value = 17975.25
$ws.Range("N71").Value = -26087.25
$ws.Range("H72").Value = 3783.182
$ws.Range("I72").Value = 3721.25
$ws.Range("J72").Value = 3796.9443
$ws.Range("K72").Value = 33491.25
$ws.Range("L72").Value = 34172.4987
$ws.Range("M72").Value = -29435.25
$ws.Range("N72").Value = -42284.4987
$ws.Range("H73").Value = 166677330
$ws.Range("J73").Value = 166677330
$ws.Range("L73").Value = 500031990
$ws.Range("N73").Value = -500034174
$ws.Range("H74").Value = 20002.4
$ws.Range("J74").Value = 20749.75
$ws.Range("L74").Value = 62249.25
$ws.Range("N74").Value = -64371.25
$ws.Range("H76").Value = 2072010.9
$ws.Range("J76").Value = 2847452.5
$ws.Range("L76").Value = 8542357.5
$ws.Range("N76").Value = -8543123.5
$ws.Range("H77").Value = 20002.4
$ws.Range("J77").Value = 20749.75
$ws.Range("L77").Value = 186747.75
$ws.Range("N77").Value = -197355.75
$ws.Range("H79").Value = 2072010.9
$ws.Range("J79").Value = 2847452.5
$ws.Range("L79").Value = 8542357.5
$ws.Range("N79").Value = -8545009.5
$ws.Range("H80").Value = 5819.654
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 5972.44
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 17917.32
$ws.Range("M80").Value = -5064
$ws.Range("N80").Value = -19789.32
$ws.Range("H81").Value = 5097
$ws.Range("I81").Value = 3291.5
$ws.Range("K81").Value = 9874.5
$ws.Range("M81").Value = -8751.5
$ws.Range("H82").Value = 5259.222
$ws.Range("J82").Value = 5500
$ws.Range("L82").Value = 16500
$ws.Range("N82").Value = -17312
$ws.Range("H83").Value = 5819.654
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 5972.44
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 53751.96
$ws.Range("M83").Value = -13320
$ws.Range("N83").Value = -63111.96
$ws.Range("H84").Value = 5097
$ws.Range("I84").Value = 3291.5
$ws.Range("K84").Value = 29623.5
$ws.Range("M84").Value = -24007.5
$ws.Range("H85").Value = 5259.222
$ws.Range("J85").Value = 5500
$ws.Range("L85").Value = 16500
$ws.Range("N85").Value = -19308
$ws.Range("H87").Value = 21166.25
$ws.Range("I87").Value = 9999
$ws.Range("J87").Value = 24888.666
$ws.Range("K87").Value = 29997
$ws.Range("L87").Value = 74665.99800000001
$ws.Range("M87").Value = -28749
$ws.Range("N87").Value = -77161.99800000001
$ws.Range("H88").Value = 8166.6665
$ws.Range("J88").Value = 8166.6665
$ws.Range("L88").Value = 24499.9995
$ws.Range("N88").Value = -25355.9995
$ws.Range("H90").Value = 21166.25
$ws.Range("I90").Value = 9999
$ws.Range("J90").Value = 24888.666
$ws.Range("K90").Value = 89991
$ws.Range("L90").Value = 223997.994
$ws.Range("M90").Value = -83751
$ws.Range("N90").Value = -236477.994
$ws.Range("H91").Value = 8166.6665
$ws.Range("J91").Value = 8166.6665
$ws.Range("L91").Value = 24499.9995
$ws.Range("N91").Value = -27463.9995
$ws.Range("H92").Value = 255
$ws.Range("I92").Value = 233.33333
$ws.Range("K92").Value = 699.99999
$ws.Range("M92").Value = 548.00001
$ws.Range("H93").Value = 3822.6365
$ws.Range("J93").Value = 6541.5
$ws.Range("L93").Value = 19624.5
$ws.Range("N93").Value = -23368.5
$ws.Range("H94").Value = 8427.571
$ws.Range("I94").Value = 2495
$ws.Range("J94").Value = 8883.923000000001
$ws.Range("K94").Value = 7485
$ws.Range("L94").Value = 26651.769
$ws.Range("M94").Value = -6809
$ws.Range("N94").Value = -28003.769
$ws.Range("H95").Value = 9620.75
$ws.Range("J95").Value = 9620.75
$ws.Range("L95").Value = 28862.25
$ws.Range("N95").Value = -32980.25
$ws.Range("H96").Value = 11333
$ws.Range("J96").Value = 11333
$ws.Range("L96").Value = 33999
$ws.Range("N96").Value = -38117
$ws.Range("H97").Value = 402
$ws.Range("I97").Value = 402
$ws.Range("K97").Value = 1206
$ws.Range("M97").Value = -710
$ws.Range("H98").Value = 1190
$ws.Range("J98").Value = 1427.4
$ws.Range("L98").Value = 4282.200000000001
$ws.Range("N98").Value = -7278.200000000001
$ws.Range("H99").Value = 5624.6665
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 6399.6
$ws.Range("K99").Value = 5250
$ws.Range("L99").Value = 19198.8
$ws.Range("M99").Value = -3004
$ws.Range("N99").Value = -23690.8
$ws.Range("H100").Value = 12400
$ws.Range("J100").Value = 14500
$ws.Range("L100").Value = 43500
$ws.Range("N100").Value = -45122
$ws.Range("H101").Value = 10999.5
$ws.Range("J101").Value = 10999.5
$ws.Range("L101").Value = 32998.5
$ws.Range("N101").Value = -37866.5
$ws.Range("H102").Value = 21966.807
$ws.Range("J102").Value = 21966.807
$ws.Range("L102").Value = 65900.421
$ws.Range("N102").Value = -70768.421
$ws.Range("H103").Value = 1075.8334
$ws.Range("I103").Value = 450
$ws.Range("J103").Value = 1388.75
$ws.Range("K103").Value = 1350
$ws.Range("L103").Value = 4166.25
$ws.Range("M103").Value = -471
$ws.Range("N103").Value = -5924.25
$ws.Range("H104").Value = 8029.9
$ws.Range("I104").Value = 2899
$ws.Range("J104").Value = 8600
$ws.Range("K104").Value = 8697
$ws.Range("L104").Value = 25800
$ws.Range("M104").Value = -6076
$ws.Range("N104").Value = -31042
$ws.Range("H105").Value = 7526.3
$ws.Range("J105").Value = 7526.3
$ws.Range("L105").Value = 22578.9
$ws.Range("N105").Value = -27820.9
$ws.Range("H106").Value = 22499.5
$ws.Range("J106").Value = 22499.5
$ws.Range("L106").Value = 67498.5
$ws.Range("N106").Value = -69390.5
$ws.Range("H136").Value = 4809.4
$ws.Range("I136").Value = 2136.5557
$ws.Range("J136").Value = 8818.666999999999
$ws.Range("K136").Value = 6409.6671
$ws.Range("L136").Value = 26456.001
$ws.Range("M136").Value = -1309.6671
$ws.Range("N136").Value = -36656.001
$ws.Range("H137").Value = 56671370
$ws.Range("I137").Value = 187500270
$ws.Range("K137").Value = 562500810
$ws.Range("M137").Value = -562495710
$ws.Range("H138").Value = 981.6667
$ws.Range("I138").Value = 981.6667
$ws.Range("K138").Value = 2945.0001
$ws.Range("M138").Value = 2194.9999
$ws.Range("H140").Value = 8702.186
$ws.Range("I140").Value = 2985.625
$ws.Range("K140").Value = 8956.875
$ws.Range("M140").Value = -3776.875
$ws.Range("H141").Value = 10812.389
$ws.Range("I141").Value = 1642.6666
$ws.Range("K141").Value = 4927.9998
$ws.Range("M141").Value = 252.0002000000004

# ---- Sheet: GSM (29 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3664.818
$ws.Range("I113").Value = 3966.6667
$ws.Range("J113").Value = 3551.625
$ws.Range("K113").Value = 3966.6667
$ws.Range("L113").Value = 3551.625
$ws.Range("M113").Value = -1796.6667
$ws.Range("N113").Value = -7891.625
$ws.Range("H122").Value = 266827.75
$ws.Range("I122").Value = 501370.1
$ws.Range("J122").Value = 8831.15
$ws.Range("K122").Value = 1504110.3
$ws.Range("L122").Value = 26493.45
$ws.Range("M122").Value = -1501660.3
$ws.Range("N122").Value = -31393.45
$ws.Range("H123").Value = 56384.8
$ws.Range("J123").Value = 56384.8
$ws.Range("L123").Value = 56384.8
$ws.Range("N123").Value = -61284.8
$ws.Range("H129").Value = 84999.5
$ws.Range("J129").Value = 84999.5
$ws.Range("L129").Value = 84999.5
$ws.Range("N129").Value = -94999.5
$ws.Range("H132").Value = 236341.86
$ws.Range("I132").Value = 3996.543
$ws.Range("J132").Value = 1252852.6
$ws.Range("K132").Value = 11989.629
$ws.Range("L132").Value = 3758557.8
$ws.Range("M132").Value = -9459.629000000001
$ws.Range("N132").Value = -3763617.8

# ---- Sheet: LTW (43 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2473
$ws.Range("I16").Value = 1727.4231
$ws.Range("J16").Value = 8934.666999999999
$ws.Range("K16").Value = 1727.4231
$ws.Range("L16").Value = 8934.666999999999
$ws.Range("M16").Value = -1557.4231
$ws.Range("N16").Value = -9274.666999999999
$ws.Range("H40").Value = 4812468
$ws.Range("I40").Value = 7357166.5
$ws.Range("J40").Value = 5815.778
$ws.Range("K40").Value = 7357166.5
$ws.Range("L40").Value = 5815.778
$ws.Range("M40").Value = -7357030.5
$ws.Range("N40").Value = -6087.778
$ws.Range("H61").Value = 1703.9286
$ws.Range("I61").Value = 1554.7778
$ws.Range("J61").Value = 1972.4
$ws.Range("K61").Value = 1554.7778
$ws.Range("L61").Value = 1972.4
$ws.Range("M61").Value = -1352.7778
$ws.Range("N61").Value = -2376.4
$ws.Range("H100").Value = 11999.833
$ws.Range("I100").Value = 11799.8
$ws.Range("K100").Value = 11799.8
$ws.Range("M100").Value = -11258.8
$ws.Range("H113").Value = 1703.9286
$ws.Range("I113").Value = 1554.7778
$ws.Range("J113").Value = 1972.4
$ws.Range("K113").Value = 1554.7778
$ws.Range("L113").Value = 1972.4
$ws.Range("M113").Value = 615.2221999999999
$ws.Range("N113").Value = -6312.4
$ws.Range("H132").Value = 3465.038
$ws.Range("I132").Value = 3206
$ws.Range("K132").Value = 9618
$ws.Range("M132").Value = -7088
$ws.Range("H136").Value = 4801.25
$ws.Range("I136").Value = 2353
$ws.Range("J136").Value = 5850.5
$ws.Range("K136").Value = 7059
$ws.Range("L136").Value = 17551.5
$ws.Range("M136").Value = -4509
$ws.Range("N136").Value = -22651.5

# ---- Sheet: WVR (41 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 12550
$ws.Range("J18").Value = 12550
$ws.Range("L18").Value = 12550
$ws.Range("N18").Value = -12896
$ws.Range("H27").Value = 37500
$ws.Range("I27").Value = 10000
$ws.Range("K27").Value = 10000
$ws.Range("M27").Value = -9931
$ws.Range("H62").Value = 17548858
$ws.Range("I62").Value = 17548858
$ws.Range("K62").Value = 17548858
$ws.Range("M62").Value = -17548234
$ws.Range("H65").Value = 17548858
$ws.Range("I65").Value = 17548858
$ws.Range("K65").Value = 87744290
$ws.Range("M65").Value = -87741170
$ws.Range("H100").Value = 524007.66
$ws.Range("I100").Value = 751140.5600000001
$ws.Range("J100").Value = 1602
$ws.Range("K100").Value = 1502281.12
$ws.Range("L100").Value = 3204
$ws.Range("M100").Value = -1501740.12
$ws.Range("N100").Value = -4286
$ws.Range("H126").Value = 3226.2856
$ws.Range("I126").Value = 2300.125
$ws.Range("K126").Value = 6900.375
$ws.Range("M126").Value = -4430.375
$ws.Range("H132").Value = 3115.7551
$ws.Range("I132").Value = 980.575
$ws.Range("J132").Value = 12605.444
$ws.Range("K132").Value = 2941.725
$ws.Range("L132").Value = 37816.33199999999
$ws.Range("M132").Value = -411.7250000000004
$ws.Range("N132").Value = -42876.33199999999
$ws.Range("H136").Value = 6859.6
$ws.Range("I136").Value = 3956.2744
$ws.Range("J136").Value = 9881.429
$ws.Range("K136").Value = 11868.8232
$ws.Range("L136").Value = 29644.287
$ws.Range("M136").Value = -9318.823199999999
$ws.Range("N136").Value = -34744.287
